# Natmi following Dr Hou advice
# Rebuild the LR-pairs output table to include the previously-missing
# "self" sending->target cluster rows (ECs->ECs, FAPs->FAPs, sCs->sCs),
# producing the full 3x3 matrix of Sending cluster x Target cluster.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$clusters = @("ECs", "FAPs", "sCs")

# data[sending][target] = array of 16 numeric values for columns E..T
$data = @{
    "ECs" = @{
        "ECs"  = @(3, 1, 22.59487733333333, 67.784632, 0.7395019553569895, 0.7395019553569895, 1, 0.3333333333333333, 0.058393, 0.175179, 0.01085373024912483, 0.01085373024912483, 1.319382672125333, 11.874444049128, 0.008026354742145116, 0.008026354742145116)
        "FAPs" = @(3, 1, 22.59487733333333, 67.784632, 0.7395019553569895, 0.7395019553569895, 3, 1, 0.4773683333333333, 1.432105, 0.08873027793527143, 0.08873027793527143, 10.78607893448444, 97.07471041036, 0.06561621403250235, 0.06561621403250235)
        "sCs"  = @(3, 1, 22.59487733333333, 67.784632, 0.7395019553569895, 0.7395019553569895, 3, 1, 4.844232333333333, 14.532697, 0.9004159918156037, 0.9004159918156038, 109.4548353458338, 985.0935181125041, 0.665859386582342, 0.6658593865823421)
    }
    "FAPs" = @{
        "ECs"  = @(3, 1, 0.3045986666666667, 0.9137960000000001, 0.00996913177602551, 0.00996913177602551, 1, 0.3333333333333333, 0.058393, 0.175179, 0.01085373024912483, 0.01085373024912483, 0.01778642994266667, 0.160077869484, 0.0001082022671149596, 0.0001082022671149596)
        "FAPs" = @(3, 1, 0.3045986666666667, 0.9137960000000001, 0.00996913177602551, 0.00996913177602551, 3, 1, 0.4773683333333333, 1.432105, 0.08873027793527143, 0.08873027793527143, 0.1454057578422222, 1.30865182058, 0.0008845638332600896, 0.0008845638332600896)
        "sCs"  = @(3, 1, 0.3045986666666667, 0.9137960000000001, 0.00996913177602551, 0.00996913177602551, 3, 1, 4.844232333333333, 14.532697, 0.9004159918156037, 0.9004159918156038, 1.475546709756889, 13.279920387812, 0.008976365675650461, 0.008976365675650463)
    }
    "sCs" = @{
        "ECs"  = @(3, 1, 7.654706, 22.964118, 0.2505289128669849, 0.2505289128669849, 1, 0.3333333333333333, 0.058393, 0.175179, 0.01085373024912483, 0.01085373024912483, 0.446981247458, 4.022831227121999, 0.002719173239864753, 0.002719173239864753)
        "FAPs" = @(3, 1, 7.654706, 22.964118, 0.2505289128669849, 0.2505289128669849, 3, 1, 0.4773683333333333, 1.432105, 0.08873027793527143, 0.08873027793527143, 3.654114245376667, 32.88702820839, 0.02222950006950897, 0.02222950006950897)
        "sCs"  = @(3, 1, 7.654706, 22.964118, 0.2505289128669849, 0.2505289128669849, 3, 1, 4.844232333333333, 14.532697, 0.9004159918156037, 0.9004159918156038, 37.08117430736067, 333.730568766246, 0.2255802395576112, 0.2255802395576112)
    }
}

$row = 2
foreach ($sending in $clusters) {
    foreach ($target in $clusters) {
        $ws.Cells.Item($row, 1).Value = $sending
        $ws.Cells.Item($row, 2).Value = "L1cam"
        $ws.Cells.Item($row, 3).Value = "Erbb3"
        $ws.Cells.Item($row, 4).Value = $target

        $values = $data[$sending][$target]
        for ($i = 0; $i -lt $values.Length; $i++) {
            $ws.Cells.Item($row, 5 + $i).Value = $values[$i]
        }

        $row++
    }
}
